# Add a new "accession_id" field to the metadata workbook.
#
# This mirrors a user:
#   1) inserting a new column B ("accession_id") on the "Metadata" sheet, and
#   2) inserting a new row 4 ("accession_id") on the "Attribute description" sheet,
# then filling in the new cells (re-using the "No" / long description text that
# already matches the style used by the other rows/cols).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Attribute description")

# --- 1) "Metadata" sheet: insert a new column B for accession_id -----------
$ws1.Columns("B:B").Insert()
$ws1.Range("B1").Value = "accession_id"

# --- 2) "Attribute description" sheet: insert a new row 4 for accession_id -
$ws2.Rows("4:4").Insert()
$ws2.Range("A4").Value = "accession_id"
$ws2.Range("B4").Value = "No"
$ws2.Range("C4").Value = "Unique identifier given to a DNA or protein sequence record to allow for tracking of different versions of that sequence record and the associated sequence over time in a single data repository (e.g. NCBI)"

# Row insert at row 4 (between the header row 3 and the former row 4) does not
# always pick up the boxed border used by the "Mandatory" column (column B) in
# the surrounding rows, so make sure it matches by copying the format down
# from the row below (row 5, which is the old row 4 "sample_name" row).
$ws2.Range("B5").Copy()
$ws2.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
